$wb = $excel.ActiveWorkbook

$wsFilepath = $wb.Worksheets.Item("Filepath")
$wsFilepath.Range("B2").Value = "../../Calculs_EcoDynBat/Generation_Data/"
$wsFilepath.Range("B3").Value = "../../Calculs_EcoDynBat/Importation_Data/"
$wsFilepath.Range("B6").Value = "../../Calculs_EcoDynBat/Mappings/Mapping_case_residue_mean.xlsx"

$wsFilepath.Activate()
